$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H5").Value = 99.25
$wsALC.Range("I5").Value = 45
$wsALC.Range("J5").Value = 189.66667
$wsALC.Range("K5").Value = 45
$wsALC.Range("L5").Value = 189.66667
$wsALC.Range("M5").Value = 70
$wsALC.Range("N5").Value = -419.66667
$wsALC.Range("H6").Value = 472.9
$wsALC.Range("I6").Value = 448.2857
$wsALC.Range("J6").Value = 530.3333
$wsALC.Range("K6").Value = 1344.8571
$wsALC.Range("L6").Value = 1590.9999
$wsALC.Range("M6").Value = -1232.8571
$wsALC.Range("N6").Value = -1814.9999
$wsALC.Range("H28").Value = 781
$wsALC.Range("I28").Value = 781
$wsALC.Range("K28").Value = 781
$wsALC.Range("M28").Value = -296
$wsALC.Range("H32").Value = 3867.4
$wsALC.Range("I32").Value = 3814
$wsALC.Range("J32").Value = 3947.5
$wsALC.Range("K32").Value = 3814
$wsALC.Range("L32").Value = 3947.5
$wsALC.Range("M32").Value = -3488
$wsALC.Range("N32").Value = -4599.5
$wsALC.Range("H41").Value = 333.58823
$wsALC.Range("I41").Value = 368.0909
$wsALC.Range("J41").Value = 270.33334
$wsALC.Range("K41").Value = 368.0909
$wsALC.Range("L41").Value = 270.33334
$wsALC.Range("M41").Value = 71.90910000000002
$wsALC.Range("N41").Value = -1150.33334
$wsALC.Range("H62").Value = 8858.352999999999
$wsALC.Range("I62").Value = 8858.352999999999
$wsALC.Range("K62").Value = 8858.352999999999
$wsALC.Range("M62").Value = -8234.352999999999
$wsALC.Range("H65").Value = 8858.352999999999
$wsALC.Range("I65").Value = 8858.352999999999
$wsALC.Range("K65").Value = 44291.765
$wsALC.Range("M65").Value = -41171.765
$wsALC.Range("H70").Value = 2948.4
$wsALC.Range("I70").Value = 0
$wsALC.Range("J70").Value = 2948.4
$wsALC.Range("K70").Value = 0
$wsALC.Range("L70").Value = 8845.200000000001
$wsALC.Range("M70").ClearContents()
$wsALC.Range("N70").Value = -9385.200000000001
$wsALC.Range("H73").Value = 2948.4
$wsALC.Range("I73").Value = 0
$wsALC.Range("J73").Value = 2948.4
$wsALC.Range("K73").Value = 0
$wsALC.Range("L73").Value = 8845.200000000001
$wsALC.Range("M73").ClearContents()
$wsALC.Range("N73").Value = -10717.2
$wsALC.Range("H137").Value = 32612.727
$wsALC.Range("I137").Value = 38956.742
$wsALC.Range("J137").Value = 4064.6667
$wsALC.Range("K137").Value = 116870.226
$wsALC.Range("L137").Value = 12194.0001
$wsALC.Range("M137").Value = -114320.226
$wsALC.Range("N137").Value = -17294.0001

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H32").Value = 28790.88
$wsARM.Range("I32").Value = 30506.027
$wsARM.Range("J32").Value = 16098.8
$wsARM.Range("K32").Value = 30506.027
$wsARM.Range("L32").Value = 16098.8
$wsARM.Range("M32").Value = -30219.027
$wsARM.Range("N32").Value = -16672.8
$wsARM.Range("H94").Value = 60000
$wsARM.Range("J94").Value = 60000
$wsARM.Range("L94").Value = 60000
$wsARM.Range("N94").Value = -61802
$wsARM.Range("H123").Value = 84000
$wsARM.Range("I123").Value = 89000
$wsARM.Range("J123").Value = 71500
$wsARM.Range("K123").Value = 89000
$wsARM.Range("L123").Value = 71500
$wsARM.Range("M123").Value = -84100
$wsARM.Range("N123").Value = -81300
$wsARM.Range("H132").Value = 32421.03
$wsARM.Range("I132").Value = 35061.355
$wsARM.Range("K132").Value = 105184.065
$wsARM.Range("M132").Value = -102654.065

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H86").Value = 3655.3333
$wsBSM.Range("I86").Value = 2899.75
$wsBSM.Range("J86").Value = 5166.5
$wsBSM.Range("K86").Value = 2899.75
$wsBSM.Range("L86").Value = 5166.5
$wsBSM.Range("M86").Value = -1776.75
$wsBSM.Range("N86").Value = -7412.5
$wsBSM.Range("H89").Value = 3655.3333
$wsBSM.Range("I89").Value = 2899.75
$wsBSM.Range("J89").Value = 5166.5
$wsBSM.Range("K89").Value = 14498.75
$wsBSM.Range("L89").Value = 25832.5
$wsBSM.Range("M89").Value = -8882.75
$wsBSM.Range("N89").Value = -37064.5
$wsBSM.Range("H105").Value = 3008.15
$wsBSM.Range("I105").Value = 2731.8333
$wsBSM.Range("K105").Value = 2731.8333
$wsBSM.Range("M105").Value = -984.8332999999998

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H7").Value = 100.26667
$wsCRP.Range("I7").Value = 39.142857
$wsCRP.Range("K7").Value = 39.142857
$wsCRP.Range("M7").Value = 73.85714300000001
$wsCRP.Range("H31").Value = 1155.5
$wsCRP.Range("I31").Value = 1155.5
$wsCRP.Range("K31").Value = 1155.5
$wsCRP.Range("M31").Value = -860.5
$wsCRP.Range("H34").Value = 1155.5
$wsCRP.Range("I34").Value = 1155.5
$wsCRP.Range("K34").Value = 1155.5
$wsCRP.Range("M34").Value = -953.5
$wsCRP.Range("H80").Value = 28127
$wsCRP.Range("J80").Value = 28127
$wsCRP.Range("L80").Value = 28127
$wsCRP.Range("N80").Value = -30373
$wsCRP.Range("H83").Value = 28127
$wsCRP.Range("J83").Value = 28127
$wsCRP.Range("L83").Value = 84381
$wsCRP.Range("N83").Value = -95613
$wsCRP.Range("H97").Value = 45049.25
$wsCRP.Range("I97").Value = 40000
$wsCRP.Range("J97").Value = 46732.332
$wsCRP.Range("K97").Value = 40000
$wsCRP.Range("L97").Value = 46732.332
$wsCRP.Range("M97").Value = -39009
$wsCRP.Range("N97").Value = -48714.332

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H57").Value = 9949.5
$wsCUL.Range("I57").Value = 8900
$wsCUL.Range("J57").Value = 10999
$wsCUL.Range("K57").Value = 26700
$wsCUL.Range("L57").Value = 32997
$wsCUL.Range("M57").Value = -26141
$wsCUL.Range("N57").Value = -34115
$wsCUL.Range("H62").Value = 10999
$wsCUL.Range("J62").Value = 10999
$wsCUL.Range("L62").Value = 32997
$wsCUL.Range("N62").Value = -34369
$wsCUL.Range("H65").Value = 10999
$wsCUL.Range("J65").Value = 10999
$wsCUL.Range("L65").Value = 98991
$wsCUL.Range("N65").Value = -105855
$wsCUL.Range("H105").Value = 10999
$wsCUL.Range("J105").Value = 10999
$wsCUL.Range("L105").Value = 32997
$wsCUL.Range("N105").Value = -38239
$wsCUL.Range("H128").Value = 189398
$wsCUL.Range("I128").Value = 189398
$wsCUL.Range("K128").Value = 568194
$wsCUL.Range("M128").Value = -563214

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H35").Value = 1748.3334
$wsLTW.Range("I35").Value = 1372.5
$wsLTW.Range("J35").Value = 2500
$wsLTW.Range("K35").Value = 1372.5
$wsLTW.Range("L35").Value = 2500
$wsLTW.Range("M35").Value = -1036.5
$wsLTW.Range("N35").Value = -3172
$wsLTW.Range("H61").Value = 1729.0769
$wsLTW.Range("I61").Value = 1865.4445
$wsLTW.Range("J61").Value = 1422.25
$wsLTW.Range("K61").Value = 1865.4445
$wsLTW.Range("L61").Value = 1422.25
$wsLTW.Range("M61").Value = -1663.4445
$wsLTW.Range("N61").Value = -1826.25
$wsLTW.Range("H80").Value = 69995.5
$wsLTW.Range("J80").Value = 69995.5
$wsLTW.Range("L80").Value = 69995.5
$wsLTW.Range("N80").Value = -72241.5
$wsLTW.Range("H83").Value = 69995.5
$wsLTW.Range("J83").Value = 69995.5
$wsLTW.Range("L83").Value = 209986.5
$wsLTW.Range("N83").Value = -221218.5
$wsLTW.Range("H113").Value = 1729.0769
$wsLTW.Range("I113").Value = 1865.4445
$wsLTW.Range("J113").Value = 1422.25
$wsLTW.Range("K113").Value = 1865.4445
$wsLTW.Range("L113").Value = 1422.25
$wsLTW.Range("M113").Value = 304.5554999999999
$wsLTW.Range("N113").Value = -5762.25
$wsLTW.Range("H122").Value = 4427.5
$wsLTW.Range("I122").Value = 3095
$wsLTW.Range("K122").Value = 9285
$wsLTW.Range("M122").Value = -6835
$wsLTW.Range("H132").Value = 70540.55499999999
$wsLTW.Range("I132").Value = 95710
$wsLTW.Range("K132").Value = 287130
$wsLTW.Range("M132").Value = -284600

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H126").Value = 54974.5
$wsWVR.Range("I126").Value = 71339.2
$wsWVR.Range("K126").Value = 214017.6
$wsWVR.Range("M126").Value = -211547.6
$wsWVR.Range("H129").Value = 0
$wsWVR.Range("J129").Value = 0
$wsWVR.Range("L129").Value = 0
$wsWVR.Range("N129").ClearContents()
